$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; remove protection before editing, then
# restore it once the cell writes are complete.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A11)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-30 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.4916335377646724
$ws.Range("E2").Value = -0.006110394459908908

$ws.Range("D3").Value = 0.2498662010856393
$ws.Range("E3").Value = -0.007420091324200739

$ws.Range("D4").Value = 0.09896549069021275
$ws.Range("E4").Value = -0.01629844259326318

$ws.Range("D5").Value = 0.1023642296615219
$ws.Range("E5").Value = -0.01130152848540988

$ws.Range("D6").Value = 0.02959980943848455
$ws.Range("E6").Value = -0.01547333012974539

$ws.Range("D7").Value = 0.02757073135946904
$ws.Range("E7").Value = -0.01195520581113796

$ws.Range("D8").Value = 1
$ws.Range("E8").Value = -0.008415581893470025

# Restore sheet protection (same semantics as before, minus the
# original obfuscated password which the COM object model cannot
# replicate bit-for-bit).
$ws.Protect()
